$wb = $excel.ActiveWorkbook

# --- Update topic_prompts sheet prompt texts (3 questions -> 2 questions) ---
$ws2 = $wb.Worksheets.Item("topic_prompts")

$b2 = @"
Assist me in drafting a 40 word clear, inspiring and memorable personal statement for myself that I can share with others.

My CliftonStrengths are: {strengths}.

Guide my thought process through a sequence of questions and proposed answers, to help identify what I should include. Please initiate the thought-provoking sequence of questions by asking me one question and only ask the next one when an answer is provided. With each question, provide several suggestions.

During the questions, help me identify my values, what's important to me, and where I want to make meaningful contributions to the world and help others.

Throughout the sequence of questions, focus on facilitating my thought-process by asking follow-up questions and offering further suggestions when appropriate. 

The questions should be asked one at a time, only move onto the next question after I have replied. 

Ask me exactly 2 questions. 

Our conversation should progress as follows:

[Brief introduction, then ask Question 1, provided by you]
%Answer 1, written by me%
[Brief acknowledgement to Answer 1 (paraphrase but don't overly repeat), then Question 2, provided by you]
%Answer 2, written by me%
...
[Brief acknowledgement to Answer 1 (paraphrase but don't overly repeat), then Question 3, provided by you]
%Answer 2, written by me%
[Brief acknowledgement to Answer 2 (paraphrase but don't overly repeat). Inform that next we will take everything that's been shared and suggest some personal statement ideas.]

Do not provide any suggestions for the personal statement at this stage, focus on understanding more about me.
"@

$b3 = @"
Assist me in drafting a 40 word clear, inspiring and memorable personal statement for myself that I can use at work, and that aligns with my company culture and values.

My Clifton Strengths are: {strengths}.

Guide my thought process through a sequence of questions and proposed answers, to help identify what I should include. Please initiate the thought-provoking sequence of questions by asking me one question and only ask the next one when an answer is provided. With each question, provide several suggestions.

During the questions, seek to understand what i do at work, help me identify my values and what's important to me, what the culture and values of my company are, and where I can make meaningful contributions.

Throughout the sequence of questions, focus on facilitating my thought-process by asking follow-up questions and offering further suggestions when appropriate. 

The questions should be asked one at a time, only move onto the next question after I have replied. 

Ask me exactly 2 questions. 

Our conversation should progress as follows:

[Brief introduction, then ask Question 1, provided by you]
%Answer 1, written by me%
[Brief acknowledgement to Answer 1 (paraphrase but don't overly repeat), then Question 2, provided by you]
%Answer 2, written by me%
...
[Brief acknowledgement to Answer 1 (paraphrase but don't overly repeat), then Question 3, provided by you]
%Answer 2, written by me%
[Brief acknowledgement to Answer 2 and share that we will next work on some personal statement ideas (paraphrase but don't overly repeat).]

Do not provide any suggestions for the personal statement at this stage, focus on understanding more about me.
"@

$b4 = @"
My Clifton Strengths are: {strengths}.

Assist me in finding ways I can form a better working relationship and get things done with a colleague, team-member or stakeholder,  leveraging our individual strengths. 

Guide my thought process through a sequence of questions and proposed answers, to help identify ways I can use my strengths to achieve desired outcomes with that individual. Please initiate the thought-provoking sequence of questions by asking me one question and only ask the next one when an answer is provided. With each question, provide several suggestions.

During the questions, seek to understand:
- more about what outcome I am hoping to achieve and what may be standing in the way 
- the strengths of my colleague and their working style (that may or may not be helpful)
- potential partnerships we can have to deliver on shared goals

Throughout the sequence of questions, focus on facilitating my thought-process by asking follow-up questions and offering further suggestions when appropriate. 

The questions should be asked one at a time, only move onto the next question after I have replied. 

Ask me exactly 2 questions. 

Our conversation should progress as follows:

[Brief introduction, then ask Question 1, provided by you]
%Answer 1, written by me%
[Brief acknowledgement to Answer 1 (paraphrase but don't overly repeat), then Question 2, provided by you]
%Answer 2, written by me%
...
[Brief reply to Answer 1, then Question 2, provided by you]
%Answer 2, written by me%
[Brief reply Answer 2 and share that we will next work on some ideas for improving collaboration (paraphrase but don't overly repeat).]

Focus on understanding my situation, any challenges in working style and the outcomes I need to achieve.
"@

$b5 = @"
My Clifton Strengths are: {strengths}.

Assist me in finding ways I can active change at work,  leveraging my strengths. 

Guide my thought process through a sequence of questions and proposed answers, to help identify ways I can use my strengths to activate change following change management principals. Please initiate the thought-provoking sequence of questions by asking me one question and only ask the next one when an answer is provided. With each question, provide several suggestions.

During the questions, seek to understand:
- more about what change I am hoping to achieve, the degree of change and it's impact on people, processes, systems, culture, behaviour, structure etc.
- any challenges I am currently or expect to face
- who are the stakeholders involved or impacted and what might their expected reactions be that could stand in the way
- the interventions currently in place or being thought of and how might one assess their effectiveness towards implementing the change

Throughout the sequence of questions, focus on facilitating my thought-process by asking follow-up questions and offering further suggestions when appropriate. 

The questions should be asked one at a time, only move onto the next question after I have replied. 

Ask me exactly 2 questions. 

Our conversation should progress as follows:

[Brief introduction, then ask Question 1, provided by you]
%Answer 1, written by me%
[Brief reply to Answer 1, then Question 2, provided by you]
%Answer 2, written by me%
...
[Brief reply to Answer 1, then Question 2, provided by you]
%Answer 2, written by me%
[Brief reply to Answer 2 and share that we will next work on some ideas for using your strengths to bring about change.]

Focus on understanding my situation and the change I need to achieve.
"@

$ws2.Range("B2").Value = $b2
$ws2.Range("B3").Value = $b3
$ws2.Range("B4").Value = $b4
$ws2.Range("B5").Value = $b5

# --- View / selection changes ---
# Active sheet becomes "topic_prompts" (was "text")
$ws2.Activate()
$ws2.Range("A3").Select()

$ws5 = $wb.Worksheets.Item("text")
$ws5.Range("B8").Select()

$ws2.Activate()
